$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-gender")

# Update is_active flag to FALSE for the "Others" (OTH) gender rows
$ws.Range("D4").Value = $false
$ws.Range("D7").Value = $false
$ws.Range("D10").Value = $false

# Update the active cell selection recorded in the sheet view
$ws.Range("D12").Select()
